$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F2 status from "In progress" to "Complete"
$ws.Range("F2").Value = "Complete"

# Add Status and Technician for row 5
$ws.Range("F5").Value = "In progress"
$ws.Range("G5").Value = "Thomas Kosacz"

# Update selected cell to G2
$ws.Range("G2").Select()
